$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.754.53"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.641.40"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.870.80"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.646.79"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.73"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "26.751.48"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.67"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.118"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.62"
$ws.Range("D29").ClearFormats()
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "1.287.87"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.539"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").Value = "1.780.78"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.406"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.25%  "
